# PRJ0018886_Hierarchy Viewer+ Time recordManager (Partial changes)
# Updates project-name text used across several sheets:
#   "Engagement Project Bend-FVA-109081" / "E - Project Clear-FVA-105379"
#     -> "Project Bend-Bernhard Capital Partners Management-FVA-109081"
#   "Project Clear-FVA-105379"
#     -> "Project Clear-LucidHealth-FVA-105379" (added as a new cell alongside)

$wb = $excel.ActiveWorkbook

$newBend  = "Project Bend-Bernhard Capital Partners Management-FVA-109081"
$newClear = "Project Clear-LucidHealth-FVA-105379"

# Sheet: Project_Title
$ws = $wb.Worksheets.Item("Project_Title")
$ws.Range("A2").Value = $newBend
$ws.Range("D2").Value = $newClear

# Sheet: RateSheetManagement
$ws = $wb.Worksheets.Item("RateSheetManagement")
$ws.Range("A2").Value = $newBend
$ws.Range("C2").Value = $newClear

# Sheet: WeeklyEntryMatrix
$ws = $wb.Worksheets.Item("WeeklyEntryMatrix")
$ws.Range("A2").Value = $newBend
$ws.Range("E2").Value = $newClear
